$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds decimal "amount_value" figures stored as text (e.g. "660.00").
# Mark the cells as Text before writing so the numeric-looking strings keep
# their literal representation (including the trailing zeros) instead of
# being coerced into plain numbers.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A10").NumberFormat = "@"

# Row 2: amount unchanged, count updated
$ws.Range("B2").Value = 117

# Row 3: amount + count updated
$ws.Range("A3").Value = "660.00"
$ws.Range("B3").Value = 6

# Row 4: amount + count updated
$ws.Range("A4").Value = "440.00"
$ws.Range("B4").Value = 17

# Row 5: amount updated, count unchanged
$ws.Range("A5").Value = "1980.00"

# Row 6: amount + count updated
$ws.Range("A6").Value = "880.00"
$ws.Range("B6").Value = 2

# Row 7: unchanged

# Row 8: amount + count updated
$ws.Range("A8").Value = "60.00"
$ws.Range("B8").Value = 2

# Row 9: amount + count updated
$ws.Range("A9").Value = "270.00"
$ws.Range("B9").Value = 2

# Row 10: amount updated, count unchanged
$ws.Range("A10").Value = "180.00"
